$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PackageTrackNum / ShipmentTrackNum values (6th April 2022 data refresh).
# These strings look numeric, so a leading apostrophe forces text storage
# (matching the original shared-string "t=s" cell type), and resetting the
# style back to Normal avoids leaving a stray NumberFormat-derived style.
$values = @{
    2  = "320018252380"
    3  = "320018252391"
    4  = "320018252428"
    5  = "320018252461"
    6  = "320018252520"
    7  = "320018252564"
    8  = "320018252612"
    9  = "320018252656"
    10 = "320018252689"
    11 = "320018252704"
    12 = "320018252748"
    13 = "320018252760"
    14 = "320018252807"
    15 = "320018252829"
    16 = "320018252873"
    17 = "320018252910"
    18 = "320018252976"
    19 = "320018253012"
    20 = "320018253240"
    21 = "320018253284"
    22 = "320018253354"
}

# Rows where column D mirrors column C (ShipmentTrackNum == PackageTrackNum).
$dMirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $values.Keys) {
    $val = $values[$row]

    $cCell = $ws.Range("C$row")
    $cCell.Value = "'" + $val
    $cCell.Style = "Normal"

    if ($dMirrorRows -contains $row) {
        $dCell = $ws.Range("D$row")
        $dCell.Value = "'" + $val
        $dCell.Style = "Normal"
    }
}
